$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores numeric-looking text (e.g. "30.075.18",
# "26.40") that must remain literal text, matching the original
# inlineStr cells. Temporarily force the column to Text format so the
# COM layer does not auto-coerce the assigned strings into numbers
# (which would also silently drop meaningful trailing zeros, e.g.
# "26.40" -> 26.4). The style is reset back to Normal afterwards so no
# stray formatting change is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.075.18'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '1.906.41'
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '0.7432'
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").Value = '243.82'
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -3.04%  '
$ws.Range("D9").Value = '26.40'
$ws.Range("E9").Value = '  -5.72%  '
$ws.Range("D10").Value = '0.06968'
$ws.Range("E10").Value = '  -4.29%  '
$ws.Range("D11").Value = '0.08084'
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("D12").Value = '0.7667'
$ws.Range("E12").Value = '  -1.75%  '
$ws.Range("D13").Value = '1.946.45'
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").Value = '5.304'
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("D15").Value = '92.15'
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").Value = '14.23'
$ws.Range("E16").Value = '  -1.88%  '
$ws.Range("D17").Value = '30.073.40'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").Value = '6.066'
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("D19").Value = '0.000007825'
$ws.Range("E19").Value = '  -2.83%  '
$ws.Range("D20").Value = '239.84'
$ws.Range("E20").Value = '  -4.84%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.155.44'
$ws.Range("E22").Value = '  -1.49%  '
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = '7.142'
$ws.Range("E24").Value = '  +6.68%  '
$ws.Range("E25").Value = '  -2.09%  '
$ws.Range("D26").Value = '166.86'
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("D27").Value = '19.00'
$ws.Range("E27").Value = '  -0.57%  '
$ws.Range("E28").Value = '  -2.69%  '
$ws.Range("D29").Value = '2.048'
$ws.Range("E29").Value = '  -6.90%  '
$ws.Range("D30").Value = '1.352'
$ws.Range("E30").Value = '  -1.63%  '
$ws.Range("D31").Value = '1.541'
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("D32").Value = '4.332'
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("D34").Value = '0.05222'
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").Value = '1.300'
$ws.Range("E35").Value = '  -2.48%  '
$ws.Range("D36").Value = '0.7469'
$ws.Range("E36").Value = '  -1.34%  '
$ws.Range("E37").Value = '  -2.52%  '
$ws.Range("D38").Value = '0.01966'
$ws.Range("E38").Value = '  +0.25%  '
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("D40").Value = '6.317'
$ws.Range("E40").Value = '  -2.83%  '
$ws.Range("D41").Value = '0.4484'
$ws.Range("E41").Value = '  -0.75%  '
$ws.Range("D42").Value = '74.23'
$ws.Range("E42").Value = '  -6.10%  '
$ws.Range("D43").Value = '1.971'
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").Value = '0.8398'
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = '7.724'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '101.66'
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("D48").Value = '9.869'
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("D49").Value = '2.074.96'
$ws.Range("E49").Value = '  -0.95%  '
$ws.Range("D50").Value = '36.63'
$ws.Range("E50").Value = '  -2.73%  '
$ws.Range("D51").Value = '0.1181'
$ws.Range("E51").Value = '  -4.44%  '

$ws.Range("D2:D51").Style = "Normal"
